$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('永辉超市', 'N沐曦-U', '平潭发展'),
    @('航天电子', '永辉超市', '中金公司'),
    @('百大集团', '平潭发展', '沐曦股份'),
    @('航天发展', '航天电子', '永辉超市'),
    @('N沐曦-U', '航天发展', '航天发展'),
    @('英维克', '百大集团', '雪人集团'),
    @('平潭发展', '摩尔线程-U', '天际股份'),
    @('航天机电', '永鼎股份', '东兴证券'),
    @('天际股份', '翠微股份', '航天电子'),
    @('雪人集团', '航天机电', '顺灏股份'),
    @('翠微股份', '雪人集团', '翠微股份'),
    @('浙江世宝', '英维克', '百大集团'),
    @('盛新锂能', '国机重装', '通宇通讯'),
    @('国机重装', '浙江世宝', '太阳电缆'),
    @('摩尔线程-U', '天际股份', '龙洲股份'),
    @('奕东电子', '华菱线缆', '浙江世宝'),
    @('永鼎股份', '盛新锂能', '信达证券'),
    @('通宇通讯', '新易盛', '航天机电'),
    @('飞龙股份', '利群股份', '英维克'),
    @('新易盛', '航天信息', '来伊份')
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}